$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "birthday" column (C2:C12) was stored as real dates (numFmtId 165,
# "yyyy-mm-dd;@"). The fix re-stores them as plain text strings formatted
# as "MM/dd/yyyy" (numFmtId 49, the builtin Text format), because the
# Discord-bot code that reads this sheet expects text, not a date serial.

# Switch the column to Text format first so the values we assign are kept
# as literal strings instead of being re-parsed back into date serials.
$ws.Range("C2:C12").NumberFormat = "@"

# Re-assign the birthday strings. (Row 6 is intentionally set using the
# same literal as row 9 -- "not assigning correctly" per the commit
# message -- reproducing the original buggy conversion.)
$ws.Range("C2").Value = "06/11/2020"
$ws.Range("C3").Value = "05/13/2020"
$ws.Range("C4").Value = "08/17/2020"
$ws.Range("C5").Value = "07/29/2020"
$ws.Range("C7").Value = "05/27/2020"
$ws.Range("C8").Value = "04/05/2020"
$ws.Range("C6").Value = "07/03/2020"
$ws.Range("C9").Value = "07/03/2020"
$ws.Range("C10").Value = "12/22/2020"
$ws.Range("C11").Value = "11/06/2020"
$ws.Range("C12").Value = "10/14/2020"

# Final selection left on F10, matching the saved workbook state.
[void]$ws.Range("F10").Select()
